$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new rows 22-30 (regcntr_id/usr_id series continuing from existing data)
$ws.Range("A22").Value = 10002
$ws.Range("B22").Value = 110021
$ws.Range("C22").Value = "eng"
$ws.Range("D22").Value = $true
$ws.Range("E22").Value = "superadmin"
$ws.Range("F22").Value = "now()"
$ws.Range("A23").Value = 10003
$ws.Range("B23").Value = 110022
$ws.Range("C23").Value = "eng"
$ws.Range("D23").Value = $true
$ws.Range("E23").Value = "superadmin"
$ws.Range("F23").Value = "now()"
$ws.Range("A24").Value = 10004
$ws.Range("B24").Value = 110023
$ws.Range("C24").Value = "eng"
$ws.Range("D24").Value = $true
$ws.Range("E24").Value = "superadmin"
$ws.Range("F24").Value = "now()"
$ws.Range("A25").Value = 10005
$ws.Range("B25").Value = 110024
$ws.Range("C25").Value = "eng"
$ws.Range("D25").Value = $true
$ws.Range("E25").Value = "superadmin"
$ws.Range("F25").Value = "now()"
$ws.Range("A26").Value = 10006
$ws.Range("B26").Value = 110025
$ws.Range("C26").Value = "eng"
$ws.Range("D26").Value = $true
$ws.Range("E26").Value = "superadmin"
$ws.Range("F26").Value = "now()"
$ws.Range("A27").Value = 10007
$ws.Range("B27").Value = 110026
$ws.Range("C27").Value = "eng"
$ws.Range("D27").Value = $true
$ws.Range("E27").Value = "superadmin"
$ws.Range("F27").Value = "now()"
$ws.Range("A28").Value = 10008
$ws.Range("B28").Value = 110027
$ws.Range("C28").Value = "eng"
$ws.Range("D28").Value = $true
$ws.Range("E28").Value = "superadmin"
$ws.Range("F28").Value = "now()"
$ws.Range("A29").Value = 10009
$ws.Range("B29").Value = 110028
$ws.Range("C29").Value = "eng"
$ws.Range("D29").Value = $true
$ws.Range("E29").Value = "superadmin"
$ws.Range("F29").Value = "now()"
$ws.Range("A30").Value = 10010
$ws.Range("B30").Value = 110029
$ws.Range("C30").Value = "eng"
$ws.Range("D30").Value = $true
$ws.Range("E30").Value = "superadmin"
$ws.Range("F30").Value = "now()"

# Select rows below the data (as left by the author after entering data)
$null = $ws.Rows("31:1048576").Select()

# Set page orientation (introduces <pageSetup> element on save)
$ws.PageSetup.Orientation = 1
